# fedexGSR.ts now retries (up to 3x) when the Fedex response has no
# Tracking Number, and denial reasons get mapped onto the output sheet.
# Re-shape the little "Tracking Number / Description" report to reflect
# that: two PENDING tracking numbers, each denied with
# "UNABLE TO PROCESS REQUEST", and drop the old sample / blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe everything first so no stale rows/columns survive outside the new
# A1:B3 used range.
$ws.Cells.Clear()

$ws.Range("A1").Value = "Tracking Number"
$ws.Range("B1").Value = "Description"

$ws.Range("A2").Value = "720266531415-PENDING"
$ws.Range("B2").Value = "UNABLE TO PROCESS REQUEST"

$ws.Range("A3").Value = "720315971720-PENDING"
$ws.Range("B3").Value = "UNABLE TO PROCESS REQUEST"

# Widen the tracking-number column and give the data rows a touch more
# height, matching the re-saved report.
$ws.Columns.Item(1).ColumnWidth = 23.75
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

# Leave the cursor on the last data row, like the saved workbook.
$ws.Range("A3").Select()
